# Minor fix in TSP.
# Update the "Fitness" column (C) values for the run's generations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-10 (Generation 0-8): Fitness 4098 -> 4500
foreach ($r in 2..10) {
    $ws.Cells.Item($r, 3).Value = 4500
}

# Rows 11-12 (Generation 9-10): Fitness 4098 -> 4567
foreach ($r in 11..12) {
    $ws.Cells.Item($r, 3).Value = 4567
}
